$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "26.776.18"; E = "  -0.11%  " },
    @{ Row = 3;  D = "1.639.83";  E = "  -0.28%  " },
    @{ Row = 4;  D = $null;       E = "  -0.07%  " },
    @{ Row = 5;  D = "218.43";    E = "  +0.72%  " },
    @{ Row = 6;  D = $null;       E = "  -0.09%  " },
    @{ Row = 7;  D = $null;       E = "  -0.21%  " },
    @{ Row = 8;  D = $null;       E = "  -0.07%  " },
    @{ Row = 9;  D = $null;       E = "  -0.84%  " },
    @{ Row = 10; D = "19.17";     E = "  +0.09%  " },
    @{ Row = 11; D = "0.0847";    E = "  +0.69%  " },
    @{ Row = 12; D = "1.869.90";  E = "  -0.15%  " },
    @{ Row = 13; D = "1.638.78";  E = "  -0.54%  " },
    @{ Row = 14; D = $null;       E = "  -0.65%  " },
    @{ Row = 15; D = $null;       E = "  -0.47%  " },
    @{ Row = 16; D = $null;       E = "  +0.76%  " },
    @{ Row = 17; D = "26.800.92"; E = "  +0.01%  " },
    @{ Row = 18; D = $null;       E = "  -0.67%  " },
    @{ Row = 19; D = "215.99";    E = "  +0.80%  " },
    @{ Row = 20; D = $null;       E = "  -0.26%  " },
    @{ Row = 21; D = $null;       E = "  -0.11%  " },
    @{ Row = 22; D = $null;       E = "  +4.51%  " },
    @{ Row = 23; D = $null;       E = "  -1.60%  " },
    @{ Row = 25; D = "147.52";    E = "  +1.73%  " },
    @{ Row = 26; D = $null;       E = "  -0.22%  " },
    @{ Row = 27; D = $null;       E = "  -0.36%  " },
    @{ Row = 28; D = "7.10";      E = "  +0.46%  " },
    @{ Row = 29; D = "15.68";     E = "  -0.04%  " },
    @{ Row = 30; D = $null;       E = "  -1.45%  " },
    @{ Row = 31; D = $null;       E = "  +0.86%  " },
    @{ Row = 32; D = $null;       E = "  +1.87%  " },
    @{ Row = 33; D = $null;       E = "  -0.87%  " },
    @{ Row = 34; D = $null;       E = "  +0.65%  " },
    @{ Row = 35; D = "1.262.46";  E = "  -2.27%  " },
    @{ Row = 36; D = $null;       E = "  +0.32%  " },
    @{ Row = 37; D = $null;       E = "  -0.25%  " },
    @{ Row = 38; D = $null;       E = "  -1.66%  " },
    @{ Row = 39; D = "0.815";     E = "  -1.31%  " },
    @{ Row = 40; D = $null;       E = "  -0.28%  " },
    @{ Row = 41; D = "0.804";     E = "  -0.78%  " },
    @{ Row = 42; D = "5.33";      E = "  -0.60%  " },
    @{ Row = 43; D = "1.780.38";  E = "  -0.80%  " },
    @{ Row = 44; D = $null;       E = "  -4.35%  " },
    @{ Row = 45; D = "92.36";     E = "  +1.00%  " },
    @{ Row = 46; D = "60.86";     E = "  +0.68%  " },
    @{ Row = 47; D = "1.59";      E = "  -0.86%  " },
    @{ Row = 48; D = $null;       E = "  -0.75%  " },
    @{ Row = 49; D = $null;       E = "  -0.71%  " },
    @{ Row = 50; D = "7.55";      E = "  -1.54%  " },
    @{ Row = 51; D = $null;       E = "  -1.86%  " }
)

$textForceRows = @(5, 10, 11, 19, 25, 28, 29, 39, 41, 42, 45, 46, 47, 50)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        if ($textForceRows -contains $u.Row) {
            # These values look numeric (e.g. "218.43"); force them to stay
            # text so they round-trip exactly like the original inline string.
            $dCell.NumberFormat = "@"
            $dCell.Value = $u.D
            $dCell.NumberFormat = "General"
        } else {
            $dCell.Value = $u.D
        }
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
